# netCrypto.xlsx update
# - Cell T2 on the active sheet changes from 425737 to 427745
# - The sheet's current selection moves from T3 to T2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in T2
$ws.Range("T2").Value = 427745

# Move the selection/active cell to T2 (was T3)
$ws.Range("T2").Select()
